$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '50.868.72'
$ws.Range("E2").Value = '  -2.58%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.901.04'
$ws.Range("E3").Value = '  -2.56%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '371.36'
$ws.Range("E5").Value = '  +4.76%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '101.49'
$ws.Range("E6").Value = '  -5.58%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.540'
$ws.Range("E7").Value = '  -4.12%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  -0.04%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.583'
$ws.Range("E9").Value = '  -5.10%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.68'
$ws.Range("E10").Value = '  -4.23%  '

$ws.Range("E11").Value = '  +0.11%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0831'
$ws.Range("E12").Value = '  -3.03%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '18.19'
$ws.Range("E13").Value = '  -5.53%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.365.21'
$ws.Range("E14").Value = '  -2.29%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.32'
$ws.Range("E15").Value = '  -4.11%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.905.91'
$ws.Range("E16").Value = '  -2.23%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.921'
$ws.Range("E17").Value = '  -7.98%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '50.881.26'
$ws.Range("E18").Value = '  -2.56%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.22'
$ws.Range("E19").Value = '  -7.96%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.17'
$ws.Range("E20").Value = '  -4.52%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.82'
$ws.Range("E21").Value = '  -5.93%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0₃0937'
$ws.Range("E22").Value = '  -3.90%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '68.01'
$ws.Range("E23").Value = '  -2.33%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '258.76'
$ws.Range("E24").Value = '  -2.02%  '

$ws.Range("E25").Value = '  -2.51%  '

$ws.Range("E26").Value = '  -6.02%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.999'
$ws.Range("E27").Value = '  -0.03%  '

$ws.Range("E28").Value = '  -4.30%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '25.51'
$ws.Range("E29").Value = '  -4.97%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.02'
$ws.Range("E30").Value = '  -7.90%  '

$ws.Range("E31").Value = '  -7.66%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.20'
$ws.Range("E32").Value = '  +1.09%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '9.81'
$ws.Range("E33").Value = '  -4.71%  '

$ws.Range("E34").Value = '  -2.87%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '51.17'
$ws.Range("E35").Value = '  +0.76%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '33.90'
$ws.Range("E36").Value = '  -7.04%  '

$ws.Range("E37").Value = '  +0.47%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0420'
$ws.Range("E38").Value = '  -5.17%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.97'
$ws.Range("E39").Value = '  -7.32%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '16.90'
$ws.Range("E40").Value = '  -5.65%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.56'
$ws.Range("E41").Value = '  -5.96%  '

$ws.Range("E42").Value = '  -7.19%  '

$ws.Range("E43").Value = '  -4.45%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '119.19'
$ws.Range("E44").Value = '  -1.98%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '21.88'
$ws.Range("E45").Value = '  -4.12%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.09'
$ws.Range("E46").Value = '  -1.77%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.009.38'
$ws.Range("E47").Value = '  -5.18%  '

$ws.Range("E48").Value = '  -1.79%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.12'
$ws.Range("E49").Value = '  -7.56%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.194.76'
$ws.Range("E50").Value = '  -2.28%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.234'
$ws.Range("E51").Value = '  -2.52%  '
